# SequentialMemoryUsage.xlsx update
# - Refresh the "Memory Usage" data table with two more rows of samples (rows 5 & 6)
#   and correct a handful of previously-zeroed/misaligned cells in rows 2-4.
# - Give each of the 4 charts an explicit title.
# - Extend each chart series' source range so the new rows are plotted.
# - Re-flow charts 3 & 4 onto a second row (2x2 grid) to make room, matching
#   the new layout; chart 1/2 positions are (almost) unchanged.
# - Restore the zoom level + selected cell recorded in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory Usage")

# ---------------------------------------------------------------------------
# 1. Worksheet data
# ---------------------------------------------------------------------------

# Row 2 (fix D2, G2, H2 which were stale/zeroed)
$ws.Range("D2").Value = 1331688
$ws.Range("G2").Value = 1547344
$ws.Range("H2").Value = 2728

# Row 3 (fix C3, D3, G3, H3)
$ws.Range("C3").Value = 1331688
$ws.Range("D3").Value = 1331688
$ws.Range("G3").Value = 1547440
$ws.Range("H3").Value = 2472

# Row 4 (fix E4, G4, H4)
$ws.Range("E4").Value = 5203096
$ws.Range("G4").Value = 1552360
$ws.Range("H4").Value = 2600

# Row 5 (new data row)
$ws.Range("A5").Value = 9176616
$ws.Range("B5").Value = 7865744
$ws.Range("C5").Value = 1331688
$ws.Range("D5").Value = 1331688
$ws.Range("E5").Value = 5203072
$ws.Range("F5").Value = 5203120
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1200

# Row 6 (new data row; C6/D6 intentionally left blank)
$ws.Range("A6").Value = 9105384
$ws.Range("B6").Value = 7804552
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 2032
$ws.Range("G6").Value = 2663560
$ws.Range("H6").Value = 0

# ---------------------------------------------------------------------------
# 2. Chart titles + series ranges
# ---------------------------------------------------------------------------

# --- Chart 1 : "Memory Usage Graph 10000"  (series A, B) ---
$chart1 = $ws.ChartObjects(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Memory Usage Graph 10000"
$chart1.SeriesCollection(1).Formula = "=SERIES('Memory Usage'!`$A`$1,,'Memory Usage'!`$A`$2:`$A`$6,1)"
$chart1.SeriesCollection(2).Formula = "=SERIES('Memory Usage'!`$B`$1,,'Memory Usage'!`$B`$2:`$B`$6,2)"

# --- Chart 2 : "Memory Usage Graph 1000"  (series C, D) ---
$chart2 = $ws.ChartObjects(2).Chart
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Memory Usage Graph 1000"
$chart2.SeriesCollection(1).Formula = "=SERIES('Memory Usage'!`$C`$1,,'Memory Usage'!`$C`$2:`$C`$5,1)"
$chart2.SeriesCollection(2).Formula = "=SERIES('Memory Usage'!`$D`$1,,'Memory Usage'!`$D`$2:`$D`$4,2)"

# --- Chart 3 : "Memory Usage Tree 10000"  (series E, F) ---
$chart3 = $ws.ChartObjects(3).Chart
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Memory Usage Tree 10000"
$chart3.SeriesCollection(1).Formula = "=SERIES('Memory Usage'!`$E`$1,,'Memory Usage'!`$E`$2:`$E`$5,1)"
$chart3.SeriesCollection(2).Formula = "=SERIES('Memory Usage'!`$F`$1,,'Memory Usage'!`$F`$2:`$F`$6,2)"

# --- Chart 4 : "Memory Usage Tree 1000"  (series G, H) ---
$chart4 = $ws.ChartObjects(4).Chart
$chart4.HasTitle = $true
$chart4.ChartTitle.Text = "Memory Usage Tree 1000"
$chart4.SeriesCollection(1).Formula = "=SERIES('Memory Usage'!`$G`$1,,'Memory Usage'!`$G`$2:`$G`$5,1)"
$chart4.SeriesCollection(2).Formula = "=SERIES('Memory Usage'!`$H`$1,,'Memory Usage'!`$H`$2:`$H`$5,2)"

# ---------------------------------------------------------------------------
# 3. Re-position charts 2 / 4 (minor nudge) and 3 / 4 (new row) into the
#    2x2 grid layout (points, matching the sheet's default col/row metrics).
# ---------------------------------------------------------------------------

$co1 = $ws.ChartObjects(1)
$co1.Left = 3
$co1.Top = 301.425
$co1.Width = 201.625
$co1.Height = 219.6

$co2 = $ws.ChartObjects(2)
$co2.Left = 213.775
$co2.Top = 299.475
$co2.Width = 124.1625
$co2.Height = 218.85

$co3 = $ws.ChartObjects(3)
$co3.Left = 3.15
$co3.Top = 524.475
$co3.Width = 201.475
$co3.Height = 218.85

$co4 = $ws.ChartObjects(4)
$co4.Left = 214.825
$co4.Top = 525.975
$co4.Width = 124.7625
$co4.Height = 218.85

# ---------------------------------------------------------------------------
# 4. Sheet view: zoom + selection
# ---------------------------------------------------------------------------

$excel.ActiveWindow.Zoom = 115
$ws.Range("J5").Select()
